# Updates the "cryptos" symbol list with refreshed price/volume figures
# (and a few re-ordered rows) per the Feb 12 2023 data refresh.
#
# Note: columns D (Price) and E (Volume(1h)) hold numeric-looking text
# (e.g. "41.00", "0.64%") that must stay stored as text, not be
# auto-converted to numbers (which would drop formatting such as
# trailing zeros or the trailing "%"). Prefixing the value with a
# leading apostrophe forces Excel to keep it as literal text, exactly
# like typing `'41.00` into a cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-0.06%"
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'0.64%"
$ws.Range("D4").Value = "'5.225"
$ws.Range("E4").Value = "'2.04%"
$ws.Range("E5").Value = "'0.67%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.311"
$ws.Range("E6").Value = "'1.36%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.634"
$ws.Range("E7").Value = "'0.87%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9153"
$ws.Range("E8").Value = "'1.60%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.437"
$ws.Range("E9").Value = "'-0.45%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1247"
$ws.Range("E10").Value = "'13.38%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1822"
$ws.Range("E11").Value = "'3.04%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09192"
$ws.Range("E12").Value = "'0.22%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04166"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1050"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001258"
$ws.Range("E15").Value = "'0.38%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005763"
$ws.Range("E16").Value = "'-0.66%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007509"
$ws.Range("E17").Value = "'2,395.62%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.344"
$ws.Range("E18").Value = "'-0.28%"
$ws.Range("D19").Value = "'0.3336"
$ws.Range("D20").Value = "'7.450"
$ws.Range("E20").Value = "'13.13%"
$ws.Range("D21").Value = "'0.1384"
$ws.Range("E21").Value = "'1.43%"
$ws.Range("D22").Value = "'0.2883"
$ws.Range("E22").Value = "'7.52%"
$ws.Range("D23").Value = "'0.04078"
$ws.Range("E23").Value = "'0.30%"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("E24").Value = "'3.32%"
$ws.Range("D25").Value = "'0.004279"
$ws.Range("E25").Value = "'4.66%"
$ws.Range("D26").Value = "'0.0001272"
$ws.Range("E26").Value = "'-2.30%"
$ws.Range("D38").Value = "'0.02490"
$ws.Range("E38").Value = "'4.89%"
$ws.Range("E39").Value = "'3.02%"
$ws.Range("D40").Value = "'0.007842"
$ws.Range("E40").Value = "'1.13%"
$ws.Range("D41").Value = "'0.1313"
$ws.Range("E41").Value = "'1.05%"
$ws.Range("D42").Value = "'0.006664"
$ws.Range("E42").Value = "'-1.36%"
$ws.Range("E43").Value = "'-2.04%"
$ws.Range("D44").Value = "'0.007662"
$ws.Range("E44").Value = "'-3.51%"
$ws.Range("E45").Value = "'-8.46%"
$ws.Range("D46").Value = "'0.00006730"
$ws.Range("E46").Value = "'-4.09%"
$ws.Range("E47").Value = "'0.15%"
$ws.Range("D48").Value = "'0.4269"
$ws.Range("E48").Value = "'1,341.42%"
$ws.Range("D49").Value = "'0.003106"
$ws.Range("E49").Value = "'-26.08%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.15%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.15%"
